# Add results workbooks: split the single "Sheet1" (parallel feeder results)
# into three scenario sheets - parallel, sectionalized, contingency - each
# holding the same PV-penetration / generation-output table but with its own
# network-losses column.

$wb = $excel.ActiveWorkbook

# --- rename the original sheet -------------------------------------------------
$parallel = $wb.Worksheets.Item(1)
$parallel.Name = "parallel"

# --- duplicate it (keeps headers, styles, column widths, row heights, etc.) ---
$parallel.Copy($null, $parallel) | Out-Null
$sectionalized = $wb.Worksheets.Item(2)
$sectionalized.Name = "sectionalized"

$sectionalized.Copy($null, $sectionalized) | Out-Null
$contingency = $wb.Worksheets.Item(3)
$contingency.Name = "contingency"

# --- update the "Network Losses (kW)" column for each new scenario ------------
$sectionalized.Range("C2").Value = 190.5
$sectionalized.Range("C3").Value = 165.9
$sectionalized.Range("C4").Value = 144
$sectionalized.Range("C5").Value = 124.7
$sectionalized.Range("C6").Value = 108
$sectionalized.Range("C7").Value = 93.9

$contingency.Range("C2").Value = 229
$contingency.Range("C3").Value = 199.7
$contingency.Range("C4").Value = 173.5
$contingency.Range("C5").Value = 150.30000000000001
$contingency.Range("C6").Value = 130.1
$contingency.Range("C7").Value = 112.9

# --- restore per-sheet selections / active tab ---------------------------------
$parallel.Cells.Select() | Out-Null
$sectionalized.Range("A1:C7").Select() | Out-Null
$contingency.Range("H23").Select() | Out-Null

# "contingency" ends up the active / visible tab
$contingency.Activate() | Out-Null
